$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 (exhibition) ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 146
$ws.Range("F3").Value = 9
$ws.Range("F4").Value = 423
$ws.Range("F5").Value = 205
$ws.Range("F7").Value = 1278
$ws.Range("F8").Value = 473
$ws.Range("F10").Value = 239
$ws.Range("F12").Value = 196
$ws.Range("F13").Value = 1076
$ws.Range("F14").Value = 10
$ws.Range("F16").Value = 13
$ws.Range("F17").Value = 222
$ws.Range("F18").Value = 1587
$ws.Range("F19").Value = 584
$ws.Range("F20").Value = 251
$ws.Range("F21").Value = 375
$ws.Range("F23").Value = 892
$ws.Range("F26").Value = 1909
$ws.Range("F27").Value = 2737
$ws.Range("F28").Value = 1512
$ws.Range("F29").Value = 75
$ws.Range("F30").Value = 76
$ws.Range("F31").Value = 530
$ws.Range("F32").Value = 832
$ws.Range("F33").Value = 1481
$ws.Range("F34").Value = 858
$ws.Range("F35").Value = 1542
$ws.Range("F36").Value = 179
$ws.Range("F38").Value = 810
$ws.Range("F39").Value = 718
$ws.Range("F40").Value = 727
$ws.Range("F41").Value = 928
$ws.Range("F42").Value = 387

# --- Sheet: 演出 (performance) ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F6").Value = 182
$ws.Range("F15").Value = 740

# --- Sheet: 全部类型 (all types) ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 146
$ws.Range("F4").Value = 423
$ws.Range("F5").Value = 205
$ws.Range("F8").Value = 182
$ws.Range("F9").Value = 1278
$ws.Range("F10").Value = 473
$ws.Range("F12").Value = 239
$ws.Range("F14").Value = 196
$ws.Range("F15").Value = 1076
$ws.Range("F16").Value = 10
$ws.Range("F19").Value = 222
$ws.Range("F20").Value = 1587
$ws.Range("F21").Value = 584
$ws.Range("F22").Value = 251
$ws.Range("F23").Value = 375
$ws.Range("F27").Value = 2737
$ws.Range("F29").Value = 1512
$ws.Range("F30").Value = 75
$ws.Range("F31").Value = 76
$ws.Range("F32").Value = 740
$ws.Range("F34").Value = 530
$ws.Range("F35").Value = 832
$ws.Range("F36").Value = 1481
$ws.Range("F39").Value = 858
$ws.Range("F40").Value = 1542
$ws.Range("F41").Value = 810
$ws.Range("F42").Value = 718
$ws.Range("F43").Value = 727
$ws.Range("F44").Value = 928
$ws.Range("F45").Value = 387
